$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.626.52'
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").Value = '1.828.91'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.15'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5344'
$ws.Range("E7").Value = '  -0.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3995'
$ws.Range("E8").Value = '  +5.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07783'
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("E10").Value = '  +2.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.01'
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.20'
$ws.Range("E12").Value = '  +3.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.329'
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.598'
$ws.Range("E14").Value = '  +2.86%  '
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").Value = '1.825.21'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.14'
$ws.Range("E17").Value = '  +3.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001094'
$ws.Range("E18").Value = '  +2.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06589'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.82'
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.097'
$ws.Range("D23").Value = '28.626.20'
$ws.Range("E23").Value = '  +1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.23'
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.234'
$ws.Range("E25").Value = '  +6.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.85'
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.68'
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").Value = '2.036.71'
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.418'
$ws.Range("E29").Value = '  +3.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.68'
$ws.Range("E30").Value = '  +2.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.155'
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.768'
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.651'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07371'
$ws.Range("E35").Value = '  +4.56%  '
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02358'
$ws.Range("E37").Value = '  +2.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.922'
$ws.Range("E38").Value = '  +5.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.220'
$ws.Range("E39").Value = '  +2.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.42'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6312'
$ws.Range("E41").Value = '  +2.15%  '
$ws.Range("E42").Value = '  +1.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.55'
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5956'
$ws.Range("E46").Value = '  +3.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.711'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.83'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.001'
$ws.Range("E49").Value = '  +3.66%  '
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06967'
$ws.Range("E51").Value = '  +2.13%  '
